# Update the "thông tin nhân viên" sheet (first sheet):
#  - R2 ("Hành chính") is replaced by "Công Nghệ Thông Tin"
#  - R3 ("Nhân sự") is replaced by "Quan Hệ Công Chúng"
# This leaves "Hành chính" / "Nhân sự" unused so they drop out of the shared
# string table, and the new "Quan Hệ Công Chúng" string gets appended.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()

$ws1.Range("R2").Value = "Công Nghệ Thông Tin"
$ws1.Range("R3").Value = "Quan Hệ Công Chúng"

# Move the active selection to R4 (matches the updated view state in the
# workbook: selection activeCell="R4" sqref="R4").
$ws1.Range("R4").Select()
